$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 (NIK) and H2 (No_HP) contain digit strings that must stay text
# (H2 even has a leading zero) -- force text formatting before writing
# so Excel doesn't auto-coerce them into numbers.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1234456278949542"

$ws.Range("B2").Value = "BG4576HI"
$ws.Range("C2").Value = "Nia Rahmadani"
$ws.Range("D2").Value = "02-08-2025 08:38"
$ws.Range("E2").Value = 60000
$ws.Range("F2").Value = "Bank Mandiri"
$ws.Range("G2").Value = "Nia Rahmadani"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "085267947261"

$ws.Range("J2").Value = "SiCepat"
